# Week 32 profile update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark C12 as submitted (0 -> 1)
$ws.Range("C12").Value = 1

# Add new profile row for swatikadu221 / Swati Singh in row 104.
# Seed the row's cell formatting from a similarly-styled existing row
# (A/C columns use the shaded "theme text" style further up the sheet)
# before writing the new values in, same as copying the row above and
# typing over it.
$ws.Range("C96").Copy()
$ws.Range("A104").PasteSpecial(-4122)
$ws.Range("C104").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("B104").Value = "swatikadu221"
$ws.Range("A104").Value = "Swati Singh"
$ws.Range("C104").Value = 0

# Scroll back to the top of the sheet and select the cell that was just
# updated (matches the new sheetView / selection saved with the workbook)
[void]$ws.Range("C12").Select()
